$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.185.53'
$ws.Range("E2").Value = '  -2.57%  '
$ws.Range("D3").Value = '1.654.25'
$ws.Range("E3").Value = '  -4.79%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = "'235.99"
$ws.Range("E5").Value = '  -2.04%  '
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = "'0.4807"
$ws.Range("E7").Value = '  -8.07%  '
$ws.Range("D8").Value = "'0.2617"
$ws.Range("E8").Value = '  -4.64%  '
$ws.Range("D9").Value = "'0.05983"
$ws.Range("E9").Value = '  -2.96%  '
$ws.Range("D10").Value = "'0.07092"
$ws.Range("E10").Value = '  -1.40%  '
$ws.Range("D11").Value = '1.667.64'
$ws.Range("E11").Value = '  -4.11%  '
$ws.Range("D12").Value = "'0.6184"
$ws.Range("E12").Value = '  -3.76%  '
$ws.Range("D13").Value = "'14.35"
$ws.Range("E13").Value = '  -3.94%  '
$ws.Range("D14").Value = "'4.578"
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").Value = "'72.97"
$ws.Range("E15").Value = '  -5.84%  '
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").Value = '25.180.52'
$ws.Range("E18").Value = '  -2.67%  '
$ws.Range("E19").Value = '  -2.78%  '
$ws.Range("D20").Value = "'0.000006528"
$ws.Range("E20").Value = '  -3.44%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = "'4.421"
$ws.Range("E21").Value = '  +3.30%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '1.886.14'
$ws.Range("E22").Value = '  -4.11%  '
$ws.Range("D23").Value = "'8.479"
$ws.Range("E23").Value = '  -1.83%  '
$ws.Range("D24").Value = "'5.262"
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("D25").Value = "'133.02"
$ws.Range("E25").Value = '  -3.95%  '
$ws.Range("D26").Value = "'14.71"
$ws.Range("E26").Value = '  -2.99%  '
$ws.Range("D27").Value = "'1.381"
$ws.Range("E27").Value = '  -8.93%  '
$ws.Range("D28").Value = "'1.704"
$ws.Range("E28").Value = '  -3.56%  '
$ws.Range("D29").Value = "'102.25"
$ws.Range("E29").Value = '  -3.19%  '
$ws.Range("D30").Value = "'3.823"
$ws.Range("E30").Value = '  -2.70%  '
$ws.Range("D31").Value = "'0.07885"
$ws.Range("E31").Value = '  -4.81%  '
$ws.Range("D32").Value = "'3.522"
$ws.Range("E32").Value = '  -4.46%  '
$ws.Range("D33").Value = "'0.04600"
$ws.Range("E33").Value = '  -0.75%  '
$ws.Range("D34").Value = "'2.617"
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("D35").Value = "'0.9385"
$ws.Range("E35").Value = '  -5.01%  '
$ws.Range("D36").Value = "'0.5839"
$ws.Range("E36").Value = '  -5.64%  '
$ws.Range("D37").Value = "'2.607"
$ws.Range("E37").Value = '  -2.67%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = "'0.8457"
$ws.Range("E38").Value = '  +14.32%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'0.01537"
$ws.Range("E39").Value = '  -4.34%  '
$ws.Range("D40").Value = "'1.003"
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("D41").Value = "'1.838"
$ws.Range("E41").Value = '  -5.39%  '
$ws.Range("D42").Value = "'98.45"
$ws.Range("D43").Value = "'0.3696"
$ws.Range("E43").Value = '  -4.01%  '
$ws.Range("D44").Value = "'4.845"
$ws.Range("E44").Value = '  -3.21%  '
$ws.Range("D45").Value = "'0.1134"
$ws.Range("E45").Value = '  +0.30%  '
$ws.Range("D46").Value = "'6.041"
$ws.Range("E46").Value = '  -3.15%  '
$ws.Range("D47").Value = "'0.05154"
$ws.Range("E47").Value = '  -1.70%  '
$ws.Range("D48").Value = "'51.97"
$ws.Range("E48").Value = '  -5.12%  '
$ws.Range("D49").Value = "'29.53"
$ws.Range("E49").Value = '  -3.21%  '
$ws.Range("D50").Value = "'1.003"
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").Value = "'7.318"
$ws.Range("E51").Value = '  -3.99%  '
